# Apply "Add data for 2022-05-22" update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the header label to reflect the new date
$ws.Name = "Through 2022-05-14"
$ws.Range("I1").Value = "2022 (through 05-14)"

# Update the May (I6) and Total (I14) figures for the 2022 column
$ws.Range("I6").Value = 53
$ws.Range("I14").Value = 605
